# Update countries & provincias Spain
# Applies the data refresh captured in the XML diff:
#  - Updated "Datos actualizados" timestamp in A1
#  - Country rows whose case counts changed, causing Hong Kong and Cabo Verde
#    (and Laos / Santa Lucia) to swap ranking positions
#  - Updated case-count figures for several countries

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1) Timestamp update in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Julio de 2020 a las 12:38"

# 2) Numeric data updates (Casos totales, Nuevos casos, Casos activos, Recuperados,
#    Casos criticos, Muertes hoy, Muertes) for rows whose figures changed.

# Row 7 - India
$ws.Range("B7").Value = 587092
$ws.Range("C7").Value = 1300
$ws.Range("D7").Value = 348487
$ws.Range("E7").Value = 221188

# Row 13 - Iran
$ws.Range("B13").Value = 230211
$ws.Range("C13").Value = 2549
$ws.Range("D13").Value = 191487
$ws.Range("E13").Value = 27766
$ws.Range("G13").Value = 141
$ws.Range("H13").Value = 10958

# Row 48 - Rumania
$ws.Range("B48").Value = 27296
$ws.Range("C48").Value = 326
$ws.Range("D48").Value = 19314
$ws.Range("E48").Value = 6315
$ws.Range("G48").Value = 16
$ws.Range("H48").Value = 1667

# Row 73 - Malasia
$ws.Range("B73").Value = 8640
$ws.Range("C73").Value = 1
$ws.Range("D73").Value = 8375
$ws.Range("E73").Value = 144

# Row 78 - Senegal
$ws.Range("B78").Value = 6925
$ws.Range("C78").Value = 132
$ws.Range("D78").Value = 4545
$ws.Range("E78").Value = 2264
$ws.Range("G78").Value = 4
$ws.Range("H78").Value = 116

# Row 91 - Bosnia y Herzegovina
$ws.Range("B91").Value = 4606
$ws.Range("C91").Value = 153
$ws.Range("D91").Value = 2432
$ws.Range("E91").Value = 1986
$ws.Range("G91").Value = 2
$ws.Range("H91").Value = 188

# Row 94 - Hungria
$ws.Range("E94").Value = 857
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 586

# Row 104 - Albania
$ws.Range("B104").Value = 2580
$ws.Range("C104").Value = 45
$ws.Range("D104").Value = 1516
$ws.Range("E104").Value = 999
$ws.Range("G104").Value = 3
$ws.Range("H104").Value = 65

# 3) Hong Kong overtakes Cabo Verde: row 125 becomes Hong Kong (with refreshed
#    figures), row 126 becomes Cabo Verde (keeping its previous figures).
$ws.Range("A125").Value = "Hong Kong"
$ws.Range("B125").Value = 1234
$ws.Range("C125").Value = 28
$ws.Range("D125").Value = 1117
$ws.Range("E125").Value = 110
$ws.Range("H125").Value = 7

$ws.Range("A126").Value = "Cabo Verde"
$ws.Range("B126").Value = 1227
$ws.Range("C126").Value = 0
$ws.Range("D126").Value = 629
$ws.Range("E126").Value = 583
$ws.Range("H126").Value = 15

# 4) Santa Lucia overtakes Laos: figures are identical for both rows, so only
#    the country labels swap places.
$ws.Range("A203").Value = "Santa Lucia"
$ws.Range("A204").Value = "Laos"
